$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.674.47"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.530.49"

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.484"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0578"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "

$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").Value = "1.746.83"
$ws.Range("E12").Value = "  -1.94%  "

$ws.Range("D13").Value = "1.527.40"
$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("E14").Value = "  -2.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.504"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "26.667.34"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "211.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("E19").Value = "  +1.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.24%  "

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.07%  "

$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "

$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("D33").Value = "1.355.05"
$ws.Range("E33").Value = "  -2.25%  "

$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("E35").Value = "  -3.58%  "

$ws.Range("E36").Value = "  -0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.934"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.796"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.86%  "

$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("E46").Value = "  -3.85%  "

$ws.Range("D47").Value = "1.662.95"
$ws.Range("E47").Value = "  -1.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.28%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0944"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
